# Fruta / hortaliza, semanal
# Update the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) columns for rows 2-28.
# Each row's set of values is replaced with the values that another row
# used to hold (a row re-shuffle caused by the weekly refresh), while all
# other columns (A, B, C, E, F, G, H, I, N, O, Q, R) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param(
        [int]$Row,
        [double]$D,
        [double]$J,
        [double]$K,
        [double]$L,
        [double]$M,
        [double]$P
    )

    $ws.Cells.Item($Row, 4).Value = $D    # D - Fecha
    $ws.Cells.Item($Row, 10).Value = $J   # J - Volumen
    $ws.Cells.Item($Row, 11).Value = $K   # K - Precio minimo
    $ws.Cells.Item($Row, 12).Value = $L   # L - Precio maximo
    $ws.Cells.Item($Row, 13).Value = $M   # M - Precio promedio ponderado
    $ws.Cells.Item($Row, 16).Value = $P   # P - Precio $/Kg
}

Set-RowData 2  44446 34 24000 25000 24500 1633
Set-RowData 3  44351 34 24000 25000 24500 1633
Set-RowData 4  44400 16 24000 25000 24500 1633
Set-RowData 5  44390 34 24000 25000 24500 1633
Set-RowData 6  44421 18 24000 25000 24500 1633
Set-RowData 7  44385 25 14000 15000 14480 965
Set-RowData 8  44397 34 23000 24000 23500 1567
Set-RowData 9  44442 28 24000 25000 24500 1633
Set-RowData 10 44411 34 25000 26000 25500 1700
Set-RowData 11 44463 25 24000 25000 24480 1632
Set-RowData 12 44435 34 24000 25000 24500 1633
Set-RowData 13 44413 25 24000 25000 24480 1632
Set-RowData 14 44460 25 24000 25000 24480 1632
Set-RowData 15 44341 36 24000 25000 24500 1633
Set-RowData 16 44336 34 24000 25000 24500 1633
Set-RowData 17 44418 16 25000 26000 25500 1700
Set-RowData 18 44329 25 23000 23000 23000 1533
Set-RowData 19 44406 25 24000 25000 24520 1635
Set-RowData 20 44343 26 23000 24000 23500 1567
Set-RowData 21 44425 25 24000 25000 24520 1635
Set-RowData 22 44453 25 25000 26000 25520 1701
Set-RowData 23 44428 16 25000 26000 25500 1700
Set-RowData 24 44349 21 24000 25000 24524 1635
Set-RowData 25 44383 25 13000 14000 13480 899
Set-RowData 26 44449 18 24000 25000 24500 1633
Set-RowData 27 44455 18 24000 25000 24500 1633
Set-RowData 28 44432 34 24000 25000 24500 1633
